# Update crypto price/volume data from latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.027.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "'2.057.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'245.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("D7").Value = "'59.28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.01%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").Value = "'0.0775"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("D12").Value = "'15.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("D13").Value = "'0.893"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.42%  "
$ws.Range("D14").Value = "'2.357.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'5.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'2.062.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'18.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "'36.984.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'74.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "'0.0₃0891"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "'5.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'238.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'10.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.99%  "
$ws.Range("D26").Value = "'169.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").Value = "'20.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  +12.72%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").Value = "'4.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.23%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("B34").Value = "'LidoDAOToken"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'2.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("B35").Value = "'BinanceUSD"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +5.73%  "
$ws.Range("D37").Value = "'0.0842"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.26%  "
$ws.Range("D38").Value = "'1.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'0.0960"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.44%  "
$ws.Range("D44").Value = "'97.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "'17.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("D46").Value = "'1.304.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("D48").Value = "'2.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'6.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("D50").Value = "'2.245.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "'44.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.30%  "

Write-Output "Updated 89 cells"
